$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F3").Value = ""
$ws.Range("K3").Value = ""

# Row 4
$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = ""
$ws.Range("K4").Value = ""

# Row 5
$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = ""
$ws.Range("K5").Value = ""

# Row 6
$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K6").Value = ""

# Row 8
$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F8").Value = ""
$ws.Range("K8").Value = ""

# Row 9
$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F9").Value = ""
$ws.Range("K9").Value = ""

# Row 10
$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F10").Value = ""
$ws.Range("K10").Value = ""

# Row 11
$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("K11").Value = ""

# Row 12
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("K12").Value = ""

# Row 13
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("F13").Value = ""
$ws.Range("K13").Value = ""

# Row 14
$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
$ws.Range("F14").Value = ""
$ws.Range("K14").Value = ""
